$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}

Replace-Text "2025-05-20 Tuesday" "2025-05-21 Wednesday"

Replace-Text "46×15=690" "14×69=966"
Replace-Text "65×88=5720" "38×36=1368"
Replace-Text "36×50=1800" "50×67=3350"
Replace-Text "17×43=731" "31×49=1519"
Replace-Text "24×24=576" "84×23=1932"

Replace-Text "89×59=5251" "86×26=2236"
Replace-Text "53×36=1908" "30×55=1650"
Replace-Text "73×20=1460" "96×97=9312"
Replace-Text "24×92=2208" "75×81=6075"
Replace-Text "50×28=1400" "75×59=4425"

Replace-Text "45×53=2385" "23×70=1610"
Replace-Text "90×93=8370" "49×24=1176"
Replace-Text "97×99=9603" "77×95=7315"
Replace-Text "23×45=1035" "58×46=2668"
Replace-Text "87×18=1566" "93×64=5952"

Replace-Text "89×67=5963" "59×85=5015"
Replace-Text "70×61=4270" "13×12=156"
Replace-Text "59×55=3245" "50×98=4900"
Replace-Text "20×82=1640" "20×11=220"
Replace-Text "68×87=5916" "39×57=2223"

Replace-Text "21×84=1764" "59×50=2950"
Replace-Text "39×64=2496" "30×19=570"
Replace-Text "29×12=348" "50×86=4300"
Replace-Text "73×83=6059" "39×87=3393"
Replace-Text "27×19=513" "46×50=2300"

Write-Output "Done"
